$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner (A1)
$ws.Range("A1").Value2 = "Datos actualizados a 9 de Agosto de 2020 a las 23:44"

# Refresh country case data (new snapshot pulled on 9 Aug 2020 23:44),
# which also re-sorts a handful of closely-ranked countries by total cases.

# Row 4
$ws.Range("B4").Value2 = 5191822
$ws.Range("C4").Value2 = 40227
$ws.Range("D4").Value2 = 2654654
$ws.Range("E4").Value2 = 2371631
$ws.Range("G4").Value2 = 447
$ws.Range("H4").Value2 = 165537

# Row 5
$ws.Range("B5").Value2 = 3035422
$ws.Range("C5").Value2 = 22053
$ws.Range("E5").Value2 = 840080
$ws.Range("G5").Value2 = 506
$ws.Range("H5").Value2 = 101049

# Row 30
$ws.Range("B30").Value2 = 95492
$ws.Range("C30").Value2 = 178
$ws.Range("D30").Value2 = 52678
$ws.Range("E30").Value2 = 37805
$ws.Range("G30").Value2 = 17
$ws.Range("H30").Value2 = 5009

# Row 38
$ws.Range("B38").Value2 = 79732
$ws.Range("C38").Value2 = 954
$ws.Range("D38").Value2 = 43744
$ws.Range("E38").Value2 = 34679
$ws.Range("G38").Value2 = 20
$ws.Range("H38").Value2 = 1309

# Row 53
$ws.Range("B53").Value2 = 44011
$ws.Range("C53").Value2 = 382
$ws.Range("D53").Value2 = 40967
$ws.Range("E53").Value2 = 2882

# Row 79
$ws.Range("E79").Value2 = 6166
$ws.Range("G79").Value2 = 1
$ws.Range("H79").Value2 = 97

# Row 91
$ws.Range("A91").Value2 = "Guinea"
$ws.Range("B91").Value2 = 7930
$ws.Range("C91").Value2 = 55
$ws.Range("D91").Value2 = 6898
$ws.Range("E91").Value2 = 982
$ws.Range("H91").Value2 = 50

# Row 92
$ws.Range("A92").Value2 = "Gabon"
$ws.Range("B92").Value2 = 7923
$ws.Range("D92").Value2 = 5704
$ws.Range("E92").Value2 = 2168
$ws.Range("H92").Value2 = 51

# Row 116
$ws.Range("A116").Value2 = "Suazilandia"
$ws.Range("B116").Value2 = 3236
$ws.Range("C116").Value2 = 108
$ws.Range("D116").Value2 = 1607
$ws.Range("E116").Value2 = 1571
$ws.Range("G116").Value2 = 2
$ws.Range("H116").Value2 = 58

# Row 117
$ws.Range("A117").Value2 = "Somalia"
$ws.Range("B117").Value2 = 3227
$ws.Range("D117").Value2 = 1728
$ws.Range("E117").Value2 = 1406
$ws.Range("H117").Value2 = 93

# Row 124
$ws.Range("B124").Value2 = 2567
$ws.Range("C124").Value2 = 2
$ws.Range("D124").Value2 = 1962

# Row 135
$ws.Range("B135").Value2 = 1916
$ws.Range("C135").Value2 = 21
$ws.Range("D135").Value2 = 1445
$ws.Range("E135").Value2 = 403

# Row 138
$ws.Range("B138").Value2 = 1672
$ws.Range("C138").Value2 = 100
$ws.Range("D138").Value2 = 567
$ws.Range("E138").Value2 = 1030
$ws.Range("G138").Value2 = 5
$ws.Range("H138").Value2 = 75

# Row 144
$ws.Range("B144").Value2 = 1237
$ws.Range("C144").Value2 = 3
$ws.Range("D144").Value2 = 723
$ws.Range("E144").Value2 = 435

# Row 145
$ws.Range("A145").Value2 = "Gambia"
$ws.Range("B145").Value2 = 1235
$ws.Range("C145").Value2 = 145
$ws.Range("D145").Value2 = 221
$ws.Range("E145").Value2 = 991
$ws.Range("G145").Value2 = 4
$ws.Range("H145").Value2 = 23

# Row 146
$ws.Range("A146").Value2 = "Republica de Chipre"
$ws.Range("B146").Value2 = 1233
$ws.Range("C146").Value2 = 0
$ws.Range("D146").Value2 = 856
$ws.Range("E146").Value2 = 358
$ws.Range("H146").Value2 = 19

# Row 147
$ws.Range("A147").Value2 = "Georgia"
$ws.Range("B147").Value2 = 1225
$ws.Range("C147").Value2 = 9
$ws.Range("D147").Value2 = 1000
$ws.Range("E147").Value2 = 208
$ws.Range("G147").Value2 = 0
$ws.Range("H147").Value2 = 17

# Row 148
$ws.Range("A148").Value2 = "Siria"
$ws.Range("B148").Value2 = 1188
$ws.Range("C148").Value2 = 63
$ws.Range("D148").Value2 = 346
$ws.Range("E148").Value2 = 790
$ws.Range("G148").Value2 = 2
$ws.Range("H148").Value2 = 52

# Row 149
$ws.Range("A149").Value2 = "Burkina Faso"
$ws.Range("B149").Value2 = 1175
$ws.Range("D149").Value2 = 974
$ws.Range("E149").Value2 = 147
$ws.Range("H149").Value2 = 54

# Row 150
$ws.Range("A150").Value2 = "Niger"
$ws.Range("B150").Value2 = 1158
$ws.Range("C150").Value2 = 1
$ws.Range("D150").Value2 = 1057
$ws.Range("E150").Value2 = 32
$ws.Range("H150").Value2 = 69

# Row 155
$ws.Range("B155").Value2 = 944
$ws.Range("C155").Value2 = 2
$ws.Range("E155").Value2 = 29

# Row 156
$ws.Range("B156").Value2 = 898
$ws.Range("C156").Value2 = 20
$ws.Range("D156").Value2 = 104
$ws.Range("E156").Value2 = 779
$ws.Range("G156").Value2 = 1
$ws.Range("H156").Value2 = 15

# Row 164
$ws.Range("A164").Value2 = "Guyana"
$ws.Range("B164").Value2 = 568
$ws.Range("C164").Value2 = 14
$ws.Range("D164").Value2 = 189
$ws.Range("E164").Value2 = 357
$ws.Range("H164").Value2 = 22

# Row 165
$ws.Range("A165").Value2 = "Aruba"
$ws.Range("B165").Value2 = 563
$ws.Range("C165").Value2 = 54
$ws.Range("D165").Value2 = 114
$ws.Range("E165").Value2 = 446
$ws.Range("H165").Value2 = 3

# Row 168
$ws.Range("B168").Value2 = 408
$ws.Range("C168").Value2 = 3
$ws.Range("D168").Value2 = 315

# Row 202
$ws.Range("A202").Value2 = "Timor Oriental"

# Row 203
$ws.Range("A203").Value2 = "Santa Lucia"
